# Weekly update: insert a new daily price record for "Feria Lagunitas de
# Puerto Montt - Piña" as row 149, pushing the existing historical rows
# (149:213) down by one to (150:214).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 149; everything below
# (rows 149-213) shifts down to 150-214, dimension grows to T214.
$ws.Rows("149:149").Insert()

# Populate the newly inserted row 149 with the new weekly record.
$ws.Range("A149").Value = 4
$ws.Range("B149").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C149").Value = "Los Lagos"
$ws.Range("D149").Value = 44636
$ws.Range("E149").Value = 10
$ws.Range("F149").Value = "Fruta"
$ws.Range("G149").Value = 100108
$ws.Range("H149").Value = "Tropicales y subtropicales"
$ws.Range("I149").Value = 100108005
$ws.Range("J149").Value = "Piña"
$ws.Range("K149").Value = "Caramelo"
$ws.Range("L149").Value = "Tercera"
$ws.Range("M149").Value = 60
$ws.Range("N149").Value = 17000
$ws.Range("O149").Value = 18000
$ws.Range("P149").Value = 17500
$ws.Range("Q149").Value = "$/caja 16 unidades"
$ws.Range("R149").Value = "Ecuador"
$ws.Range("S149").Value = 1094
$ws.Range("T149").Value = 16
